$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.032.93'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -4.01%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.483.12'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.94%  '

$ws.Range("E4").Value = '  +0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.66'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.27%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.14'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.00%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.481.53'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.96%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.485'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.52%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.142'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.83%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.49'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.00%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.426'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.64%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000214'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -6.53%  '

$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.080.30'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.70%  '

$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '31.72'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.22%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.504.30'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.42%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.258.28'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.75%  '

$ws.Range("E18").Value = '  -0.71%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.45'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.66%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.34'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.44%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.89'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.98%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '443.19'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.60%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.621'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.12%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.83'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.14%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.634.93'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.61%  '

$ws.Range("E26").Value = '  -0.09%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.65'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.89%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000120'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -11.36%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.87'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -8.13%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.43'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -9.15%  '

$ws.Range("E31").Value = '  -5.47%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.62'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.02%  '

$ws.Range("B33").Value = 'Binance-PegBSC-USD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.08%  '

$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.166'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.35%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '25.46'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.34%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.15'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.19%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.484.40'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.81%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.81'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -8.10%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.92'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.92%  '

$ws.Range("E40").Value = '  +0.00%  '

$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.10%  '

$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.23'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -8.69%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '174.03'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.22%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0891'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.52%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.39'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.63%  '

$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.894'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.48%  '

$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '30.07'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.90%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '46.39'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.11%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.27'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -7.92%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.53'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.99%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.45'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -11.69%  '
